{"js": "// Replace the 100 equation answers in the single 20x5 table, cell by cell,\n// in row-major order. Several \"before\" equation strings repeat at different\n// positions (e.g. \"17+54=71\" and \"67+0=67\" each occur twice with different\n// replacements), so the edit must be positional (row/col), not a text\n// search-and-replace.\nconst newValues = [\n  [\"19+39=58\", \"53+14=67\", \"48-3=45\", \"31+27=58\", \"39+21=60\"],\n  [\"49+27=76\", \"66-38=28\", \"39-20=19\", \"14+17=31\", \"11+38=49\"],\n  [\"51-14=37\", \"76-17=59\", \"28+70=98\", \"55-9=46\", \"40+12=52\"],\n  [\"81+7=88\", \"63+5=68\", \"96-41=55\", \"43+53=96\", \"47-12=35\"],\n  [\"9+33=42\", \"39+12=51\", \"40+8=48\", \"86-15=71\", \"52+45=97\"],\n  [\"99-68=31\", \"20+73=93\", \"22+20=42\", \"74-72=2\", \"81-46=35\"],\n  [\"32+67=99\", \"16+21=37\", \"42-39=3\", \"73+8=81\", \"23+30=53\"],\n  [\"25+34=59\", \"71-57=14\", \"53-39=14\", \"76-60=16\", \"9+51=60\"],\n  [\"37+50=87\", \"20-4=16\", \"41-41=0\", \"25+70=95\", \"33-10=23\"],\n  [\"47+8=55\", \"97-79=18\", \"64-14=50\", \"14+37=51\", \"22-12=10\"],\n  [\"54-35=19\", \"68+15=83\", \"75-52=23\", \"57+6=63\", \"74-2=72\"],\n  [\"84+14=98\", \"91-71=20\", \"63+19=82\", \"58-51=7\", \"73-30=43\"],\n  [\"73-46=27\", \"41+42=83\", \"55-6=49\", \"27+58=85\", \"82-24=58\"],\n  [\"31+24=55\", \"22+5=27\", \"24-20=4\", \"57-41=16\", \"97-90=7\"],\n  [\"71+4=75\", \"9+39=48\", \"40+25=65\", \"32+44=76\", \"18+76=94\"],\n  [\"91-73=18\", \"21+19=40\", \"34+14=48\", \"87-51=36\", \"29+45=74\"],\n  [\"80-2=78\", \"44+0=44\", \"20+39=59\", \"29-3=26\", \"68-5=63\"],\n  [\"24-13=11\", \"88-12=76\", \"87-37=50\", \"38+20=58\", \"43-34=9\"],\n  [\"65-37=28\", \"77-23=54\", \"98-20=78\", \"67-66=1\", \"4+8=12\"],\n  [\"61-60=1\", \"61-37=24\", \"62-31=31\", \"38+11=49\", \"15+7=22\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,columnCount,values\");\nawait context.sync();\n\n// First pass: figure out which cells actually change, grab their first\n// paragraph (load queued, not yet available) for each.\nconst targets = [];\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    const oldText = table.values[r][c];\n    const newText = newValues[r][c];\n    if (oldText === newText) continue;\n    const cell = table.getCell(r, c);\n    cell.body.paragraphs.load(\"items\");\n    targets.push({ cell, newText });\n  }\n}\nawait context.sync();\n\n// Second pass: rewrite just the paragraph's text (insertText(..., \"Replace\")\n// on a Paragraph maps to a Range.Text assignment), which preserves the\n// run's existing rFonts/sz formatting and the paragraph's alignment --\n// only the <w:t> content changes, matching the diff.\nfor (const { cell, newText } of targets) {\n  const para = cell.body.paragraphs.items[0];\n  para.insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Replace the 100 equation answers in the single 20x5 table, cell by cell,\n# in row-major order. Several \"before\" equation strings repeat at different\n# positions (e.g. \"17+54=71\" and \"67+0=67\" each occur twice with different\n# replacements), so the edit must be positional (row/col), not a text\n# search-and-replace.\n$newValues = @(\n    @(\"19+39=58\", \"53+14=67\", \"48-3=45\", \"31+27=58\", \"39+21=60\"),\n    @(\"49+27=76\", \"66-38=28\", \"39-20=19\", \"14+17=31\", \"11+38=49\"),\n    @(\"51-14=37\", \"76-17=59\", \"28+70=98\", \"55-9=46\", \"40+12=52\"),\n    @(\"81+7=88\", \"63+5=68\", \"96-41=55\", \"43+53=96\", \"47-12=35\"),\n    @(\"9+33=42\", \"39+12=51\", \"40+8=48\", \"86-15=71\", \"52+45=97\"),\n    @(\"99-68=31\", \"20+73=93\", \"22+20=42\", \"74-72=2\", \"81-46=35\"),\n    @(\"32+67=99\", \"16+21=37\", \"42-39=3\", \"73+8=81\", \"23+30=53\"),\n    @(\"25+34=59\", \"71-57=14\", \"53-39=14\", \"76-60=16\", \"9+51=60\"),\n    @(\"37+50=87\", \"20-4=16\", \"41-41=0\", \"25+70=95\", \"33-10=23\"),\n    @(\"47+8=55\", \"97-79=18\", \"64-14=50\", \"14+37=51\", \"22-12=10\"),\n    @(\"54-35=19\", \"68+15=83\", \"75-52=23\", \"57+6=63\", \"74-2=72\"),\n    @(\"84+14=98\", \"91-71=20\", \"63+19=82\", \"58-51=7\", \"73-30=43\"),\n    @(\"73-46=27\", \"41+42=83\", \"55-6=49\", \"27+58=85\", \"82-24=58\"),\n    @(\"31+24=55\", \"22+5=27\", \"24-20=4\", \"57-41=16\", \"97-90=7\"),\n    @(\"71+4=75\", \"9+39=48\", \"40+25=65\", \"32+44=76\", \"18+76=94\"),\n    @(\"91-73=18\", \"21+19=40\", \"34+14=48\", \"87-51=36\", \"29+45=74\"),\n    @(\"80-2=78\", \"44+0=44\", \"20+39=59\", \"29-3=26\", \"68-5=63\"),\n    @(\"24-13=11\", \"88-12=76\", \"87-37=50\", \"38+20=58\", \"43-34=9\"),\n    @(\"65-37=28\", \"77-23=54\", \"98-20=78\", \"67-66=1\", \"4+8=12\"),\n    @(\"61-60=1\", \"61-37=24\", \"62-31=31\", \"38+11=49\", \"15+7=22\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $newText = $newValues[$r - 1][$c - 1]\n        $cell = $t.Cell($r, $c)\n\n        # Cell.Range.Text includes the trailing paragraph-mark + cell-mark\n        # characters (CR + cell-delimiter), so strip the last two characters\n        # before comparing against the plain answer text.\n        $rawText = $cell.Range.Text\n        $curText = $rawText.Substring(0, $rawText.Length - 2)\n\n        if ($curText -ne $newText) {\n            # Assigning Range.Text rewrites just the <w:t> content of the\n            # existing run in place (Word re-uses the run's current\n            # formatting and automatically keeps the paragraph/cell marks),\n            # leaving rFonts/sz and paragraph alignment untouched -- matching\n            # the diff exactly.\n            $cell.Range.Text = $newText\n        }\n    }\n}\n"}
